$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3326.5
$ws.Range("J43").Value = 3318.8333
$ws.Range("L43").Value = 3318.8333
$ws.Range("N43").Value = -3456.8333
$ws.Range("H106").Value = 8916.5
$ws.Range("I106").Value = 8770.286
$ws.Range("K106").Value = 8770.286
$ws.Range("M106").Value = -8139.286
$ws.Range("H116").Value = 20747.5
$ws.Range("I116").Value = 8993.333000000001
$ws.Range("J116").Value = 27800
$ws.Range("K116").Value = 8993.333000000001
$ws.Range("L116").Value = 27800
$ws.Range("M116").Value = -5551.333000000001
$ws.Range("N116").Value = -34684
$ws.Range("H129").Value = 3758.5
$ws.Range("J129").Value = 8841
$ws.Range("L129").Value = 26523
$ws.Range("N129").Value = -36523
$ws.Range("H138").Value = 2832.0967
$ws.Range("I138").Value = 1159.3448
$ws.Range("J138").Value = 4302.091
$ws.Range("K138").Value = 3478.0344
$ws.Range("L138").Value = 12906.273
$ws.Range("M138").Value = 1661.9656
$ws.Range("N138").Value = -23186.273
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2037.5
$ws.Range("I4").Value = 2037.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2037.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -1921.5
$ws.Range("H61").Value = 5128371
$ws.Range("I61").Value = 5720752
$ws.Range("K61").Value = 5720752
$ws.Range("M61").Value = -5720540
$ws.Range("H132").Value = 2443719.8
$ws.Range("I132").Value = 4935
$ws.Range("K132").Value = 14805
$ws.Range("M132").Value = -12275
$ws.Range("H136").Value = 5128371
$ws.Range("I136").Value = 5720752
$ws.Range("K136").Value = 17162256
$ws.Range("M136").Value = -17159706
$ws.Range("H139").Value = 139799.8
$ws.Range("I139").Value = 19000
$ws.Range("J139").Value = 169999.75
$ws.Range("K139").Value = 19000
$ws.Range("L139").Value = 169999.75
$ws.Range("M139").Value = -13860
$ws.Range("N139").Value = -180279.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 600
$ws.Range("J16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -940
$ws.Range("H20").Value = 1464.6552
$ws.Range("I20").Value = 1248.7368
$ws.Range("J20").Value = 1874.9
$ws.Range("K20").Value = 1248.7368
$ws.Range("L20").Value = 1874.9
$ws.Range("M20").Value = -1001.7368
$ws.Range("N20").Value = -2368.9
$ws.Range("H22").Value = 4999.5
$ws.Range("I22").Value = 4999.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4999.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -4826.5
$ws.Range("H86").Value = 2069.3
$ws.Range("I86").Value = 1517.1875
$ws.Range("J86").Value = 4277.75
$ws.Range("K86").Value = 1517.1875
$ws.Range("L86").Value = 4277.75
$ws.Range("M86").Value = -394.1875
$ws.Range("N86").Value = -6523.75
$ws.Range("H89").Value = 2069.3
$ws.Range("I89").Value = 1517.1875
$ws.Range("J89").Value = 4277.75
$ws.Range("K89").Value = 7585.9375
$ws.Range("L89").Value = 21388.75
$ws.Range("M89").Value = -1969.9375
$ws.Range("N89").Value = -32620.75
$ws.Range("H99").Value = 2453
$ws.Range("I99").Value = 2453
$ws.Range("K99").Value = 2453
$ws.Range("M99").Value = -955
$ws.Range("H134").Value = 25002224
$ws.Range("I134").Value = 2900
$ws.Range("J134").Value = 33335332
$ws.Range("K134").Value = 8700
$ws.Range("L134").Value = 100005996
$ws.Range("M134").Value = -6165
$ws.Range("N134").Value = -100011066
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33673104
$ws.Range("I31").Value = 47621948
$ws.Range("J31").Value = 1125801.5
$ws.Range("K31").Value = 47621948
$ws.Range("L31").Value = 1125801.5
$ws.Range("M31").Value = -47621653
$ws.Range("N31").Value = -1126391.5
$ws.Range("H34").Value = 33673104
$ws.Range("I34").Value = 47621948
$ws.Range("J34").Value = 1125801.5
$ws.Range("K34").Value = 47621948
$ws.Range("L34").Value = 1125801.5
$ws.Range("M34").Value = -47621746
$ws.Range("N34").Value = -1126205.5
$ws.Range("H58").Value = 3480.2856
$ws.Range("I58").Value = 3409.6667
$ws.Range("J58").Value = 3607.4
$ws.Range("K58").Value = 3409.6667
$ws.Range("L58").Value = 3607.4
$ws.Range("M58").Value = -3206.6667
$ws.Range("N58").Value = -4013.4
$ws.Range("H99").Value = 17246.572
$ws.Range("I99").Value = 9095.375
$ws.Range("K99").Value = 9095.375
$ws.Range("M99").Value = -7597.375
$ws.Range("H126").Value = 17246.572
$ws.Range("I126").Value = 9095.375
$ws.Range("K126").Value = 27286.125
$ws.Range("M126").Value = -24816.125
$ws.Range("H134").Value = 4088.8572
$ws.Range("I134").Value = 4103.75
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 12311.25
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -9776.25
$ws.Range("N134").Value = -17068.5
$ws.Range("H136").Value = 3480.2856
$ws.Range("I136").Value = 3409.6667
$ws.Range("J136").Value = 3607.4
$ws.Range("K136").Value = 10229.0001
$ws.Range("L136").Value = 10822.2
$ws.Range("M136").Value = -7679.000100000001
$ws.Range("N136").Value = -15922.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 18899.428
$ws.Range("I56").Value = 18899.428
$ws.Range("K56").Value = 18899.428
$ws.Range("M56").Value = -18369.428
$ws.Range("H76").Value = 10569
$ws.Range("I76").Value = 7127.6
$ws.Range("K76").Value = 21382.8
$ws.Range("M76").Value = -20999.8
$ws.Range("H79").Value = 10569
$ws.Range("I79").Value = 7127.6
$ws.Range("K79").Value = 21382.8
$ws.Range("M79").Value = -20056.8
$ws.Range("H129").Value = 5371.706
$ws.Range("I129").Value = 3730.3333
$ws.Range("K129").Value = 11190.9999
$ws.Range("M129").Value = -6190.999899999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 6799.6
$ws.Range("J10").Value = 6799.6
$ws.Range("L10").Value = 6799.6
$ws.Range("N10").Value = -7137.6
$ws.Range("H14").Value = 4749.5
$ws.Range("I14").Value = 7999
$ws.Range("K14").Value = 7999
$ws.Range("M14").Value = -7831
$ws.Range("H70").Value = 11055.117
$ws.Range("I70").Value = 8279.666999999999
$ws.Range("K70").Value = 8279.666999999999
$ws.Range("M70").Value = -8009.666999999999
$ws.Range("H73").Value = 11055.117
$ws.Range("I73").Value = 8279.666999999999
$ws.Range("K73").Value = 8279.666999999999
$ws.Range("M73").Value = -7343.666999999999
$ws.Range("H122").Value = 4539.16
$ws.Range("I122").Value = 4603.2915
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 13809.8745
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -11359.8745
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 10049837
$ws.Range("I132").Value = 2144.875
$ws.Range("K132").Value = 6434.625
$ws.Range("M132").Value = -3904.625
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("M135").Value = -160140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1738.4706
$ws.Range("I16").Value = 664.7742
$ws.Range("J16").Value = 12833.333
$ws.Range("K16").Value = 664.7742
$ws.Range("L16").Value = 12833.333
$ws.Range("M16").Value = -494.7742
$ws.Range("N16").Value = -13173.333
$ws.Range("H21").Value = 1323.3334
$ws.Range("I21").Value = 1600
$ws.Range("J21").Value = 770
$ws.Range("K21").Value = 1600
$ws.Range("L21").Value = 770
$ws.Range("M21").Value = -1426
$ws.Range("N21").Value = -1118
$ws.Range("H136").Value = 3121.04
$ws.Range("I136").Value = 2267.4
$ws.Range("K136").Value = 6802.200000000001
$ws.Range("M136").Value = -4252.200000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""
$ws.Range("H21").Value = 16417.8
$ws.Range("I21").Value = 13272.5
$ws.Range("K21").Value = 13272.5
$ws.Range("M21").Value = -13037.5
$ws.Range("H24").Value = 9000
$ws.Range("I24").Value = 9000
$ws.Range("K24").Value = 9000
$ws.Range("M24").Value = -8770
$ws.Range("H31").Value = 15333.333
$ws.Range("J31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("N31").Value = -30696
$ws.Range("H35").Value = 16417.8
$ws.Range("I35").Value = 13272.5
$ws.Range("K35").Value = 13272.5
$ws.Range("M35").Value = -12982.5
$ws.Range("H107").Value = 2491.366
$ws.Range("I107").Value = 1460.6522
$ws.Range("K107").Value = 4381.9566
$ws.Range("M107").Value = -2461.9566
$ws.Range("H126").Value = 7797.8
$ws.Range("I126").Value = 7797.8
$ws.Range("K126").Value = 23393.4
$ws.Range("M126").Value = -20923.4
$ws.Range("H132").Value = 358364.2
$ws.Range("I132").Value = 1338.3182
$ws.Range("J132").Value = 1667459
$ws.Range("K132").Value = 4014.9546
$ws.Range("L132").Value = 5002377
$ws.Range("M132").Value = -1484.9546
$ws.Range("N132").Value = -5007437
$ws.Range("H136").Value = 192069.45
$ws.Range("I136").Value = 3425.1372
$ws.Range("K136").Value = 10275.4116
$ws.Range("M136").Value = -7725.411599999999

Write-Host "Applied all changes"